# Delete the row for SSA (Salvador, Brazil) at row 214.
# This shifts all subsequent rows up by one, matching the target diff
# (dimension shrinks from A1:G330 to A1:G329).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(214).Delete()
